$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lithuania A Lyga")

# ---------------------------------------------------------------
# 1) Swap the (non-id) data of rows 89 and 90.
#    Column A (the running "id" number) must stay where it is;
#    everything from column B through AB moves between the rows.
# ---------------------------------------------------------------
$row89 = $ws.Range("B89:AB89").Value2
$row90 = $ws.Range("B90:AB90").Value2
$ws.Range("B89:AB89").Value2 = $row90
$ws.Range("B90:AB90").Value2 = $row89

# ---------------------------------------------------------------
# 2) Swap the (non-id) data of rows 101 and 102.
# ---------------------------------------------------------------
$row101 = $ws.Range("B101:AB101").Value2
$row102 = $ws.Range("B102:AB102").Value2
$ws.Range("B101:AB101").Value2 = $row102
$ws.Range("B102:AB102").Value2 = $row101

# ---------------------------------------------------------------
# 3) Rotate the (non-id) data among rows 100, 103 and 104:
#       new100 = old104 ; new103 = old100 ; new104 = old103
# ---------------------------------------------------------------
$row100 = $ws.Range("B100:AB100").Value2
$row103 = $ws.Range("B103:AB103").Value2
$row104 = $ws.Range("B104:AB104").Value2
$ws.Range("B100:AB100").Value2 = $row104
$ws.Range("B103:AB103").Value2 = $row100
$ws.Range("B104:AB104").Value2 = $row103

# ---------------------------------------------------------------
# 4) Append 5 new match rows (158-162) at the bottom of the sheet.
# ---------------------------------------------------------------
# Row 158
    $ws.Range("A158").Value = 156
    $ws.Range("B158").Value = 7862937
    $ws.Range("C158").Value = "Lithuania A Lyga"
    $ws.Range("D158").Value = 45413.375
    $ws.Range("E158").Value = "FK Siauliai"
    $ws.Range("F158").Value = "Hegelmann Litauen"
    $ws.Range("G158").Value = 1
    $ws.Range("H158").Value = 1
    $ws.Range("I158").Value = "D"
    $ws.Range("J158").Value = 2.75
    $ws.Range("K158").Value = 3.05
    $ws.Range("L158").Value = 2.375
    $ws.Range("M158").Value = 3.1
    $ws.Range("N158").Value = 3.25
    $ws.Range("O158").Value = 2.05
    $ws.Range("P158").Value = 0.25
    $ws.Range("Q158").Value = 1.95
    $ws.Range("R158").Value = 1.85
    $ws.Range("S158").Value = 2.25
    $ws.Range("T158").Value = 1.825
    $ws.Range("U158").Value = 1.975
    $ws.Range("V158").Value = -1
    $ws.Range("W158").Value = 2.25
    $ws.Range("X158").Value = -1
    $ws.Range("Y158").Value = 0.475
    $ws.Range("Z158").Value = -0.5
    $ws.Range("AA158").Value = -0.5
    $ws.Range("AB158").Value = 0.4875

# Row 159
    $ws.Range("A159").Value = 157
    $ws.Range("B159").Value = 7862052
    $ws.Range("C159").Value = "Lithuania A Lyga"
    $ws.Range("D159").Value = 45413.45833333334
    $ws.Range("E159").Value = "FK Zalgiris Vilnius"
    $ws.Range("F159").Value = "Banga Gargzdai"
    $ws.Range("G159").Value = 2
    $ws.Range("H159").Value = 2
    $ws.Range("I159").Value = "D"
    $ws.Range("J159").Value = 1.25
    $ws.Range("K159").Value = 4.8
    $ws.Range("L159").Value = 9.5
    $ws.Range("M159").Value = 1.25
    $ws.Range("N159").Value = 5
    $ws.Range("O159").Value = 8.5
    $ws.Range("P159").Value = -1.5
    $ws.Range("Q159").Value = 1.8
    $ws.Range("R159").Value = 2
    $ws.Range("S159").Value = 2.75
    $ws.Range("T159").Value = 1.825
    $ws.Range("U159").Value = 1.975
    $ws.Range("V159").Value = -1
    $ws.Range("W159").Value = 4
    $ws.Range("X159").Value = -1
    $ws.Range("Y159").Value = -1
    $ws.Range("Z159").Value = 1
    $ws.Range("AA159").Value = 0.825
    $ws.Range("AB159").Value = -1

# Row 160
    $ws.Range("A160").Value = 158
    $ws.Range("B160").Value = 7862939
    $ws.Range("C160").Value = "Lithuania A Lyga"
    $ws.Range("D160").Value = 45414.5
    $ws.Range("E160").Value = "FK Transinvest"
    $ws.Range("F160").Value = "FK Kauno Zalgiris"
    $ws.Range("G160").Value = 3
    $ws.Range("H160").Value = 2
    $ws.Range("I160").Value = "H"
    $ws.Range("J160").Value = 4.1
    $ws.Range("K160").Value = 3.3
    $ws.Range("L160").Value = 1.8
    $ws.Range("M160").Value = 5
    $ws.Range("N160").Value = 3.75
    $ws.Range("O160").Value = 1.571
    $ws.Range("P160").Value = 1
    $ws.Range("Q160").Value = 1.8
    $ws.Range("R160").Value = 2
    $ws.Range("S160").Value = 2.5
    $ws.Range("T160").Value = 2
    $ws.Range("U160").Value = 1.8
    $ws.Range("V160").Value = 4
    $ws.Range("W160").Value = -1
    $ws.Range("X160").Value = -1
    $ws.Range("Y160").Value = 0.8
    $ws.Range("Z160").Value = -1
    $ws.Range("AA160").Value = 1
    $ws.Range("AB160").Value = -1

# Row 161
    $ws.Range("A161").Value = 159
    $ws.Range("B161").Value = 7862938
    $ws.Range("C161").Value = "Lithuania A Lyga"
    $ws.Range("D161").Value = 45414.5
    $ws.Range("E161").Value = "FK Dziugas Telsiai"
    $ws.Range("F161").Value = "FK Dainava Alytus"
    $ws.Range("G161").Value = 2
    $ws.Range("H161").Value = 1
    $ws.Range("I161").Value = "H"
    $ws.Range("J161").Value = 2.35
    $ws.Range("K161").Value = 2.9
    $ws.Range("L161").Value = 3
    $ws.Range("M161").Value = 2.55
    $ws.Range("N161").Value = 2.75
    $ws.Range("O161").Value = 2.875
    $ws.Range("P161").Value = 0
    $ws.Range("Q161").Value = 1.825
    $ws.Range("R161").Value = 1.975
    $ws.Range("S161").Value = 1.75
    $ws.Range("T161").Value = 1.825
    $ws.Range("U161").Value = 1.975
    $ws.Range("V161").Value = 1.55
    $ws.Range("W161").Value = -1
    $ws.Range("X161").Value = -1
    $ws.Range("Y161").Value = 0.825
    $ws.Range("Z161").Value = -1
    $ws.Range("AA161").Value = 0.825
    $ws.Range("AB161").Value = -1

# Row 162
    $ws.Range("A162").Value = 160
    $ws.Range("B162").Value = 7862053
    $ws.Range("C162").Value = "Lithuania A Lyga"
    $ws.Range("D162").Value = 45414.54166666666
    $ws.Range("E162").Value = "Panevezys"
    $ws.Range("F162").Value = "Suduva Marijampole"
    $ws.Range("G162").Value = 1
    $ws.Range("H162").Value = 0
    $ws.Range("I162").Value = "H"
    $ws.Range("J162").Value = 1.869
    $ws.Range("K162").Value = 3.1
    $ws.Range("L162").Value = 4.1
    $ws.Range("M162").Value = 1.75
    $ws.Range("N162").Value = 3.1
    $ws.Range("O162").Value = 5
    $ws.Range("P162").Value = -0.5
    $ws.Range("Q162").Value = 1.8
    $ws.Range("R162").Value = 2
    $ws.Range("S162").Value = 2
    $ws.Range("T162").Value = 2
    $ws.Range("U162").Value = 1.8
    $ws.Range("V162").Value = 0.75
    $ws.Range("W162").Value = -1
    $ws.Range("X162").Value = -1
    $ws.Range("Y162").Value = 0.8
    $ws.Range("Z162").Value = -1
    $ws.Range("AA162").Value = -1
    $ws.Range("AB162").Value = 0.8

# ---------------------------------------------------------------
# 5) Apply the same cell formatting used by the previous last row
#    (row 157) to the 5 new rows - this carries over the bold /
#    centered / bordered style used for column A ("id") and the
#    custom date-time number format used for column D ("Date").
# ---------------------------------------------------------------
$ws.Range("A157:AB157").Copy() | Out-Null
$ws.Range("A158:AB162").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
